$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 411
$ws.Range("F6").Value = 461
$ws.Range("F7").Value = 7473
$ws.Range("F8").Value = 90
$ws.Range("F9").Value = 104
$ws.Range("F10").Value = 2074
$ws.Range("F11").Value = 8127
$ws.Range("F14").Value = 5560
$ws.Range("F16").Value = 2523
$ws.Range("F17").Value = 1074
$ws.Range("F20").Value = 395
$ws.Range("F21").Value = 86
$ws.Range("F23").Value = 446
$ws.Range("F24").Value = 1500
$ws.Range("F25").Value = 25
$ws.Range("F26").Value = 2602
$ws.Range("F28").Value = 301
$ws.Range("F29").Value = 104
$ws.Range("F30").Value = 230
$ws.Range("F33").Value = 303
$ws.Range("F34").Value = 1586
$ws.Range("F35").Value = 38
$ws.Range("F37").Value = 2516
$ws.Range("F41").Value = 304

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 92
$ws.Range("F3").Value = 92

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 411
$ws.Range("F7").Value = 461
$ws.Range("F8").Value = 7473
$ws.Range("F9").Value = 90
$ws.Range("F10").Value = 104
$ws.Range("F11").Value = 2074
$ws.Range("F12").Value = 8127
$ws.Range("F15").Value = 5560
$ws.Range("F17").Value = 2523
$ws.Range("F18").Value = 1074
$ws.Range("F21").Value = 395
$ws.Range("F22").Value = 86
$ws.Range("F23").Value = 92
$ws.Range("F25").Value = 92
$ws.Range("F26").Value = 446
$ws.Range("F27").Value = 1500
$ws.Range("F28").Value = 25
$ws.Range("F29").Value = 2602
$ws.Range("F31").Value = 301
$ws.Range("F32").Value = 104
$ws.Range("F33").Value = 230
$ws.Range("F37").Value = 303
$ws.Range("F39").Value = 1586
$ws.Range("F40").Value = 38
$ws.Range("F42").Value = 2516
$ws.Range("F47").Value = 304
